$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Wrap the institution / thesis / supervisor cell text in LaTeX \href{}{}
#    commands (commit: "More LaTeX code in excel files").
# ---------------------------------------------------------------------------

# Row 2 - PhD at University of Stirling
$ws.Range("C2").Value = '\href{https://www.stir.ac.uk/}{University of Stirling}'
$ws.Range("E2").Value = '\href{https://dspace.stir.ac.uk/handle/1893/21102}{\textbf{\textit{Contextual musicality: vocal modulation and its perception in human social interaction}}}'
$ws.Range("E2").Value = 'Tesis: ' + $ws.Range("E2").Text

$ws.Range("E3").Value = 'Supervisores: \href{https://www.scraigroberts.com/}{Prof. S. Craig Roberts}, y \href{https://scholar.google.com/citations?user=iDDoxVsAAAAJ}{Prof. Anthony C. Little}'
$ws.Range("E4").Value = 'Miembros del comité: \href{https://scholar.google.co.uk/citations?user=wxh9svQAAAAJ}{Prof. Phyllis C. Lee} (dissertation chair), y \href{https://scholar.google.com/citations?user=Qo23OGoAAAAJ}{Prof. Stuart Semple}'

# Row 5 - MSc at University of Liverpool
$ws.Range("C5").Value = '\href{https://www.liverpool.ac.uk/}{University of Liverpool}'
$ws.Range("E5").Value = 'Supervisor: \href{https://www.scraigroberts.com/}{Prof. S. Craig Roberts}'

# Row 7 - Licenciatura at Universidad Pedagogica Nacional
$ws.Range("C7").Value = '\href{https://www.upn.edu.co/}{Universidad Pedagógica Nacional}'

Write-Host "content updated"

# ---------------------------------------------------------------------------
# 2. New cell format: left/top aligned, wrapped text, applied to the whole
#    used range A1:E7. Format A1 directly then copy/paste the format onto
#    the rest of the range so every cell lands on the SAME new style index.
# ---------------------------------------------------------------------------

$ws.Range("A1").HorizontalAlignment = -4131
$ws.Range("A1").VerticalAlignment = -4160
$ws.Range("A1").WrapText = $true

$ws.Range("A1").Copy()
$ws.Range("A1:E7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Host "formatting applied"

# ---------------------------------------------------------------------------
# 3. Row heights (grown to fit the new wrapped, multi-line text).
# ---------------------------------------------------------------------------

$ws.Rows.Item(2).RowHeight = 60
$ws.Rows.Item(3).RowHeight = 60
$ws.Rows.Item(4).RowHeight = 75
$ws.Rows.Item(5).RowHeight = 45
$ws.Rows.Item(7).RowHeight = 45

Write-Host "rows resized"

# ---------------------------------------------------------------------------
# 4. Column widths.
# ---------------------------------------------------------------------------

$ws.Columns.Item(1).ColumnWidth = 30.666666666666664
$ws.Columns.Item(2).ColumnWidth = 5.166666666666666
$ws.Columns.Item(3).ColumnWidth = 21.666666666666664
$ws.Columns.Item(4).ColumnWidth = 18.666666666666664
$ws.Columns.Item(5).ColumnWidth = 57.66666666666667

Write-Host "columns resized"

# ---------------------------------------------------------------------------
# 5. Selection moves from C14 to E7.
# ---------------------------------------------------------------------------

$ws.Range("E7").Select()

Write-Host "done"
